# Update the most recent login timestamp recorded in the "Student Details"
# sheet. The scheduler app appends a new timestamp each time a user logs in;
# here the latest logged timestamp (B2) advances to the newest entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "03/14/2020 14:29:10"
